# "cash purchases for projects"
#
# The Project Budgets sheet is restructured: the old "revised budget" /
# "committed open/received/invoiced" breakdown columns (I:Q) are replaced
# with a smaller "committed budget" / "received budget" breakdown (4
# columns), and the sample row's unit of measure changes from "Hours" to
# "Piece".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old revised_budgeted_*/committed_open_*/committed_received_*/
# committed_invoiced_* columns (I through Q).
$ws.Range("I1:Q1").EntireColumn.Delete()

# Make room for the 4 replacement columns right before actual_quantity /
# actual_amount (now sitting in I:J after the delete above).
$ws.Range("I1:L1").EntireColumn.Insert()

# Sample row now uses "Piece" instead of "Hours" as its unit of measure.
$ws.Range("E2").Value = "Piece"

# New column headers.
$ws.Range("I1").Value = "committed_budget_quantity"
$ws.Range("J1").Value = "committed_budget_amount"
$ws.Range("K1").Value = "received_budget_quantity"
$ws.Range("L1").Value = "received_budget_amount"

# New column values for the sample row.
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
